# Insert a new data row at row 64 (pushing existing rows 64..151 down to 65..152),
# and populate the new row with the new price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(64).Insert()

$ws.Cells.Item(64, 1).Value = 4
$ws.Cells.Item(64, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(64, 3).Value = "Los Lagos"
$ws.Cells.Item(64, 4).Value = 44482
$ws.Cells.Item(64, 5).Value = 10
$ws.Cells.Item(64, 6).Value = "Fruta"
$ws.Cells.Item(64, 7).Value = 100102
$ws.Cells.Item(64, 8).Value = "Cítricos"
$ws.Cells.Item(64, 9).Value = 100102006
$ws.Cells.Item(64, 10).Value = "Pomelo"
$ws.Cells.Item(64, 11).Value = "Start Ruby"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 80
$ws.Cells.Item(64, 14).Value = 11000
$ws.Cells.Item(64, 15).Value = 12000
$ws.Cells.Item(64, 16).Value = 11500
$ws.Cells.Item(64, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(64, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(64, 19).Value = 821
$ws.Cells.Item(64, 20).Value = 14
